# TC33_Canine_Filter_Breed-Labrador.xlsx
# Reformat the "StatQuery" (column C) Cypher query on the startup sheet
# from a single long line into an indented, multi-line, human-readable
# version - same query, same tabs (CasesTab/SamplesTab/FilesTab rows),
# just reformatted text - and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = "MATCH (s:study)`n" +
            "  WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies`n" +
            "  MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies`n" +
            "  MATCH (d:diagnosis)`n" +
            "  WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies`n" +
            "  MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n" +
            "    WHERE demo.breed IN ['Labrador Retriever']`n" +
            "  OPTIONAL MATCH (f:file)-[*]->(c)`n" +
            "  OPTIONAL MATCH (samp:sample)-[*]->(c)`n" +
            "  WITH DISTINCT c AS c, p, s, demo, diag, f, samp`n" +
            "  RETURN count(DISTINCT(f)) as number_of_files ,`n" +
            "             count(DISTINCT(samp)) as number_of_sample ,`n" +
            "             count(DISTINCT(c.case_id)) as number_of_cases ,`n" +
            "             count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Rows 2-4 (CasesTab / SamplesTab / FilesTab) all share the same StatQuery text.
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Move the active selection from C2 to B2, scrolled so row 2 is at the top.
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
